$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Modelo" header in F1, matching the style used by existing headers (E1)
$ws.Range("F1").Value = "Modelo"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the numeric metrics in row 2
$ws.Range("B2").Value = 0.07669419371074328
$ws.Range("C2").Value = 0.9992725212023039
$ws.Range("D2").Value = 0.207546847481027

# Add the new model name value in F2
$ws.Range("F2").Value = "Pipeline(steps=[('model', AdaBoostRegressor())])"
